$d = $word.ActiveDocument

# The document currently has a single empty paragraph. Put the cursor
# at the very start of the document and type the new heading text.
$r = $d.Paragraphs(1).Range
$r.Text = "Introduction"

# Format the run we just typed (Calibri, bold, underlined, 12pt).
$r = $d.Paragraphs(1).Range
$r.Font.Name = "Calibri"
$r.Font.Bold = $true
$r.Font.Underline = 1
$r.Font.Size = 12

# Insert a new paragraph after it for the following (empty) paragraph.
$endOfFirst = $d.Paragraphs(1).Range.End
$r2 = $d.Range($endOfFirst, $endOfFirst)
$r2.InsertParagraphAfter()

# Format the new, second paragraph: Calibri 12pt, no bold/underline.
$p2 = $d.Paragraphs(2).Range
$p2.Font.Name = "Calibri"
$p2.Font.Bold = $false
$p2.Font.Underline = 0
$p2.Font.Size = 12
